$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.126.32'
$ws.Range("E2").Value = '  -1.15%  '

$ws.Range("D3").Value = '1.803.01'
$ws.Range("E3").Value = '  +0.36%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.31'
$ws.Range("E5").Value = '  +1.11%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.02%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5339'
$ws.Range("E7").Value = '  -0.11%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3784'
$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07488'
$ws.Range("E9").Value = '  -0.61%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.94'
$ws.Range("E10").Value = '  -1.68%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.099'
$ws.Range("E11").Value = '  -1.44%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  +0.02%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.223'
$ws.Range("E13").Value = '  +0.73%  '

$ws.Range("B14").Value = 'Solana'
$ws.Range("C14").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.58'
$ws.Range("E14").Value = '  -1.96%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.359'
$ws.Range("E15").Value = '  -0.19%  '

$ws.Range("D16").Value = '1.807.48'
$ws.Range("E16").Value = '  +1.07%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '89.45'
$ws.Range("E17").Value = '  -1.09%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001066'
$ws.Range("E18").Value = '  -0.11%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06515'
$ws.Range("E19").Value = '  +1.14%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9998'
$ws.Range("E20").Value = '  -0.09%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.20'
$ws.Range("E21").Value = '  -0.33%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.914'
$ws.Range("E22").Value = '  -0.11%  '

$ws.Range("D23").Value = '28.162.83'
$ws.Range("E23").Value = '  -0.95%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.15'
$ws.Range("E24").Value = '  -1.91%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.086'
$ws.Range("E25").Value = '  -2.27%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.40'
$ws.Range("E26").Value = '  -2.81%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.51'
$ws.Range("E27").Value = '  -0.17%  '

$ws.Range("D28").Value = '2.013.26'
$ws.Range("E28").Value = '  +0.87%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.313'
$ws.Range("E29").Value = '  -3.50%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '122.11'
$ws.Range("E30").Value = '  -1.07%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.125'
$ws.Range("E31").Value = '  +0.73%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1089'
$ws.Range("E32").Value = '  +7.54%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.584'
$ws.Range("E33").Value = '  -2.23%  '

$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.605'
$ws.Range("E34").Value = '  -1.71%  '

$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07161'
$ws.Range("E35").Value = '  +8.96%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2227'
$ws.Range("E36").Value = '  -3.68%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02302'
$ws.Range("E37").Value = '  -0.86%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.065'
$ws.Range("E38").Value = '  -0.78%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.538'
$ws.Range("E39").Value = '  -1.54%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6172'
$ws.Range("E40").Value = '  -2.47%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.12'
$ws.Range("E41").Value = '  -3.73%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.185'
$ws.Range("E42").Value = '  -1.91%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.430'
$ws.Range("E43").Value = '  +1.92%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.000'
$ws.Range("E44").Value = '  -0.01%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.40'
$ws.Range("E45").Value = '  -1.29%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.679'
$ws.Range("E46").Value = '  +0.16%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5765'
$ws.Range("E47").Value = '  -3.16%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.51'
$ws.Range("E48").Value = '  +0.18%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.193'
$ws.Range("E49").Value = '  +1.88%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.925'
$ws.Range("E50").Value = '  -3.24%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06822'
$ws.Range("E51").Value = '  -1.62%  '

